# Updated performance matrix list 04/20/2021
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sam")

# ---------------------------------------------------------------------
# 1) Write all string-valued cells in the exact order needed so the
#    shared-string table is built up with the same index assignment as
#    the target workbook (index 0..2 already exist: Model, Accuracy,
#    Parameters(Optimizer, Learning Rate, Epochs)).
# ---------------------------------------------------------------------
$ws.Range("A3").Value = "ResNet50"            # sst 3
$ws.Range("E2").Value = "Imagesize"           # sst 4
$ws.Range("F2").Value = "Centercrop"          # sst 5
$ws.Range("G2").Value = "normalization"       # sst 6
$ws.Range("H2").Value = "Optimizer"           # sst 7
$ws.Range("I2").Value = "Learning Rate"       # sst 8
$ws.Range("B1").Value = "Dataset"             # sst 9
$ws.Range("B2").Value = "Train"               # sst 10
$ws.Range("C2").Value = "Validation"          # sst 11
$ws.Range("D2").Value = "Test"                # sst 12
$ws.Range("G3").Value = "No"                  # sst 13
$ws.Range("H3").Value = "SGD"                 # sst 14
$ws.Range("E1").Value = "Parameters"          # sst 15
$ws.Range("K2").Value = "Accuracy "           # sst 16
$ws.Range("L2").Value = "Precision"           # sst 17
$ws.Range("M2").Value = " Recall"             # sst 18
$ws.Range("N2").Value = " F1-score"           # sst 19
$ws.Range("K1").Value = "Validation Measures" # sst 20
$ws.Range("A5").Value = "ResNet101"           # sst 21
$ws.Range("J2").Value = "Epoch"               # sst 22

# Re-used string values (no new shared-string entries)
$ws.Range("A4").Value = "ResNet50"
$ws.Range("G4").Value = "No"
$ws.Range("G5").Value = "No"
$ws.Range("H4").Value = "SGD"
$ws.Range("H5").Value = "SGD"

Write-Output "strings done"

# ---------------------------------------------------------------------
# 2) Numeric data for the three result rows.
# ---------------------------------------------------------------------
# Row 3 - ResNet50 (Centercrop=200)
$ws.Range("B3").Value = 70
$ws.Range("C3").Value = 15
$ws.Range("D3").Value = 15
$ws.Range("E3").Value = 224
$ws.Range("F3").Value = 200
$ws.Range("I3").Value = 0.01
$ws.Range("J3").Value = 25
$ws.Range("K3").Value = 0.93146399999999996
$ws.Range("L3").Value = 0.91048899999999999
$ws.Range("M3").Value = 0.91822700000000002
$ws.Range("N3").Value = 0.91416900000000001

# Row 4 - ResNet50 (Centercrop=224)
$ws.Range("B4").Value = 70
$ws.Range("C4").Value = 15
$ws.Range("D4").Value = 15
$ws.Range("E4").Value = 224
$ws.Range("F4").Value = 224
$ws.Range("I4").Value = 0.01
$ws.Range("J4").Value = 25
$ws.Range("K4").Value = 0.93561799999999995
$ws.Range("L4").Value = 0.91192300000000004
$ws.Range("M4").Value = 0.93480399999999997
$ws.Range("N4").Value = 0.92291800000000002

# Row 5 - ResNet101
$ws.Range("B5").Value = 70
$ws.Range("C5").Value = 15
$ws.Range("D5").Value = 15
$ws.Range("E5").Value = 224
$ws.Range("F5").Value = 224
$ws.Range("I5").Value = 0.01
$ws.Range("J5").Value = 25
$ws.Range("K5").Value = 0.93665600000000004
$ws.Range("L5").Value = 0.94180600000000003
$ws.Range("M5").Value = 0.91849800000000004
$ws.Range("N5").Value = 0.92884699999999998

Write-Output "numbers done"

# ---------------------------------------------------------------------
# 3) Column widths (COM ColumnWidth applies a font-dependent padding, so
#    the inputs below are pre-compensated to land on the nearest
#    achievable stored width to the target).
# ---------------------------------------------------------------------
$ws.Range("A1:D1").ColumnWidth = 19.669999999999796   # -> stored 20.5
$ws.Range("E1").ColumnWidth = 18.829999999999814       # -> stored ~19.6640625
$ws.Range("F1").ColumnWidth = 16.66999999999986        # -> stored 17.5
$ws.Range("G1").ColumnWidth = 12.82999999999994        # -> stored ~13.6640625
$ws.Range("H1").ColumnWidth = 14.6699999999999         # -> stored 15.5
$ws.Range("I1").ColumnWidth = 19.499999999999797       # -> stored ~20.33203125

Write-Output "columns done"

# ---------------------------------------------------------------------
# 4) Merge the header bands.
# ---------------------------------------------------------------------
$ws.Range("B1:D1").Merge()
$ws.Range("E1:J1").Merge()
$ws.Range("A1:A2").Merge()
$ws.Range("K1:N1").Merge()

Write-Output "merges done"

# ---------------------------------------------------------------------
# 5) Borders: every populated cell in A1:N5 gets a thin box border.
# ---------------------------------------------------------------------
$ws.Range("A1:N5").Borders.LineStyle = 1

Write-Output "borders done"

# ---------------------------------------------------------------------
# 6) Alignment.
# ---------------------------------------------------------------------
# Row 1 header band: centered + wrap text everywhere.
$ws.Range("A1:N1").HorizontalAlignment = -4108   # xlCenter
$ws.Range("A1:N1").WrapText = $true
# A1:A2 also vertically centered.
$ws.Range("A1:A2").VerticalAlignment = -4108     # xlCenter

# Row 2 sub-headers: centered only (no wrap).
$ws.Range("A2:N2").HorizontalAlignment = -4108   # xlCenter

# Model name column in data rows: left aligned.
$ws.Range("A3:A5").HorizontalAlignment = -4131   # xlLeft

Write-Output "alignment done"

# ---------------------------------------------------------------------
# 7) Fonts (bold header cells; merges that have a non-bold continuation
#    need the bold applied to the anchor cell only).
# ---------------------------------------------------------------------
$ws.Range("B1").Font.Bold = $true
$ws.Range("E1").Font.Bold = $true
$ws.Range("K1:N1").Font.Bold = $true
$ws.Range("A2:N2").Font.Bold = $true
$ws.Range("A3:A5").Font.Bold = $true

Write-Output "fonts done"

# ---------------------------------------------------------------------
# 8) Number format for the learning-rate column.
# ---------------------------------------------------------------------
$ws.Range("I3:I5").NumberFormat = "0.E+00"

Write-Output "numfmt done"

# ---------------------------------------------------------------------
# 9) Decorative bottom cell (K6): bigger Courier New font, no border,
#    taller row.
# ---------------------------------------------------------------------
$ws.Range("K6").Font.Name = "Courier New"
$ws.Range("K6").Font.Size = 14
$ws.Range("K6").Font.Bold = $false
$ws.Range("K6").Font.Color = 2171169
$ws.Rows.Item(6).RowHeight = 19

Write-Output "k6 done"

# ---------------------------------------------------------------------
# 10) Selection.
# ---------------------------------------------------------------------
$ws.Range("L9").Select()

Write-Output "done"
